$p = $ppt.ActivePresentation

# --- Slide 1 (Title slide): "TITLE" -> "Stencil Pattern" ---
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Stencil Pattern"

# --- Slide 1 (Title slide): merge the "Department / of Computer and
#     Information / Science" runs of paragraph 3 into one run, leaving
#     paragraphs 1 ("Parallel Computing") and 2 ("CIS 410/510") intact.
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange
$deptRange = $subtitleRange.Characters(32, 46)
$deptRange.Text = "Department of Computer and Information Science"

# --- Slide 2 (Overview slide): "Overview" -> "Table of Contents" ---
$slide2 = $p.Slides.Item(2)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "Table of Contents"

# --- New slide 4: duplicate slide 3 (same "Title and Content" layout,
#     with footer + slide number placeholders already populated), then
#     retitle it and clear its content placeholder. ---
$slide3 = $p.Slides.Item(3)
$dup = $slide3.Duplicate()
$slide4 = $dup.Item(1)
$slide4.Shapes.Item(1).TextFrame.TextRange.Text = "Example Implementation"
$slide4.Shapes.Item(2).TextFrame.TextRange.Text = ""
